# Update cryptos list prices/volumes (Sat Jun 17 14:06:29 UTC 2023 snapshot).
# Numeric-looking price strings are prefixed with a leading single quote so
# Excel stores them as text (matching the source data, which uses dotted
# "thousand"-style strings like "26.532.69" that are not real numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.532.69'
$ws.Range("E2").Value = '  +4.06%  '
$ws.Range("D3").Value = '1.738.00'
$ws.Range("E3").Value = '  +4.44%  '
$ws.Range("D4").Value = '''0.9999'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''244.81'
$ws.Range("E5").Value = '  +4.02%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '''0.4801'
$ws.Range("E7").Value = '  +3.79%  '
$ws.Range("D8").Value = '''0.2673'
$ws.Range("E8").Value = '  +3.96%  '
$ws.Range("D9").Value = '''0.06231'
$ws.Range("E9").Value = '  +1.83%  '
$ws.Range("D10").Value = '1.738.31'
$ws.Range("E10").Value = '  +4.55%  '
$ws.Range("D11").Value = '''0.07139'
$ws.Range("E11").Value = '  +2.68%  '
$ws.Range("D12").Value = '''15.79'
$ws.Range("E12").Value = '  +8.13%  '
$ws.Range("D13").Value = '''0.6185'
$ws.Range("E13").Value = '  +8.09%  '
$ws.Range("D14").Value = '''4.541'
$ws.Range("E14").Value = '  +4.60%  '
$ws.Range("D15").Value = '''76.97'
$ws.Range("E15").Value = '  +2.61%  '
$ws.Range("D16").Value = '''1.000'
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").Value = '26.540.19'
$ws.Range("E17").Value = '  +4.13%  '
$ws.Range("D18").Value = '''1.000'
$ws.Range("D19").Value = '''0.000006892'
$ws.Range("E19").Value = '  +2.88%  '
$ws.Range("E20").Value = '  +3.39%  '
$ws.Range("D21").Value = '1.961.98'
$ws.Range("E21").Value = '  +4.44%  '
$ws.Range("D22").Value = '''4.571'
$ws.Range("E22").Value = '  +3.70%  '
$ws.Range("D23").Value = '''8.899'
$ws.Range("E23").Value = '  +2.75%  '
$ws.Range("D24").Value = '''5.340'
$ws.Range("E24").Value = '  +2.41%  '
$ws.Range("D25").Value = '''135.52'
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("E26").Value = '  +3.59%  '
$ws.Range("D27").Value = '''1.805'
$ws.Range("E27").Value = '  +5.72%  '
$ws.Range("D28").Value = '''1.416'
$ws.Range("E28").Value = '  +3.71%  '
$ws.Range("D29").Value = '''106.93'
$ws.Range("E29").Value = '  +2.83%  '
$ws.Range("D30").Value = '''3.983'
$ws.Range("E30").Value = '  +0.75%  '
$ws.Range("D31").Value = '''3.726'
$ws.Range("E31").Value = '  +3.84%  '
$ws.Range("D32").Value = '''0.07894'
$ws.Range("E32").Value = '  +2.49%  '
$ws.Range("D33").Value = '''0.04577'
$ws.Range("E33").Value = '  +5.84%  '
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("D35").Value = '''0.9977'
$ws.Range("E35").Value = '  +6.07%  '
$ws.Range("D36").Value = '''0.6349'
$ws.Range("E36").Value = '  +6.20%  '
$ws.Range("D37").Value = '''0.9277'
$ws.Range("E37").Value = '  +1.58%  '
$ws.Range("D38").Value = '''111.66'
$ws.Range("E38").Value = '  +5.45%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '''2.436'
$ws.Range("E39").Value = '  -1.33%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '''1.983'
$ws.Range("E40").Value = '  +8.20%  '
$ws.Range("D41").Value = '''1.003'
$ws.Range("E41").Value = '  +0.37%  '
$ws.Range("E42").Value = '  +4.00%  '
$ws.Range("D43").Value = '''5.701'
$ws.Range("E43").Value = '  +14.75%  '
$ws.Range("D44").Value = '''0.3904'
$ws.Range("E44").Value = '  +5.36%  '
$ws.Range("D45").Value = '''6.927'
$ws.Range("E45").Value = '  +13.33%  '
$ws.Range("D46").Value = '''0.1196'
$ws.Range("E46").Value = '  +8.09%  '
$ws.Range("D47").Value = '''0.05331'
$ws.Range("E47").Value = '  +1.51%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''7.935'
$ws.Range("E48").Value = '  +5.69%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '''30.81'
$ws.Range("E49").Value = '  +2.92%  '
$ws.Range("D50").Value = '''1.253'
$ws.Range("E50").Value = '  +5.79%  '
$ws.Range("D51").Value = '''0.3441'
$ws.Range("E51").Value = '  +4.35%  '
